$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H104:J107").Borders.LineStyle = -4142
$ws.Range("G105").Value = "vsrinivasan2203@altimetrik.com"
